# Insert a new data row at row 116 (pushing the existing rows 116-174 down
# to 117-175) and populate it with the new weekly price record for
# Berenjena / Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("116:116").Insert()

$ws.Cells.Item(116, 1).Value2  = 3
$ws.Cells.Item(116, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(116, 3).Value2  = "Coquimbo"
$ws.Cells.Item(116, 4).Value2  = 44510
$ws.Cells.Item(116, 5).Value2  = 5
$ws.Cells.Item(116, 6).Value2  = 100112001
$ws.Cells.Item(116, 7).Value2  = "Berenjena"
$ws.Cells.Item(116, 8).Value2  = "Sin especificar"
$ws.Cells.Item(116, 9).Value2  = "Primera"
$ws.Cells.Item(116, 10).Value2 = 85
$ws.Cells.Item(116, 11).Value2 = 7500
$ws.Cells.Item(116, 12).Value2 = 8000
$ws.Cells.Item(116, 13).Value2 = 7735
$ws.Cells.Item(116, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(116, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(116, 16).Value2 = 129
$ws.Cells.Item(116, 17).Value2 = 60
$ws.Cells.Item(116, 18).Value2 = "Hortaliza"
